# feat: add 2022-Q4 data
#
# 1. Insert a new "2022-Q4" worksheet (copied from "2022-Q3" so it inherits
#    the same sheetPr / pageMargins / header styling), positioned right after
#    "总计" and before "2022-Q3". Fill it with the Q4 fund-holdings table.
# 2. Update the "总计" (summary) sheet so it gets a new first data row for
#    2022-Q4, pushing the existing 2022-Q3 / 2022-Q2 rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q4" sheet by duplicating "2022-Q3" (this keeps
# sheetPr/pageMargins/header formatting identical), then strip it down to
# just the header row and overwrite it with the Q4 data.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)

$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Source sheet had 9 data rows (rows 2-10); Q4 only needs 6 (rows 2-7).
$q4.Rows("8:10").Delete()
$q4.Range("A2:H7").ClearContents()

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Columns B, D, E, F, G hold numeric-looking text (fund code / percentages);
# force Text format up front so Excel doesn't coerce them into numbers.
$q4.Range("B2:B7").NumberFormat = "@"
$q4.Range("D2:G7").NumberFormat = "@"

$q4Data = @(
    @(0, "010054", "万家健康产业混合A",         "5.08", "91.05", "3.34", "0.1697", 8),
    @(1, "010055", "万家健康产业混合C",         "3.85", "91.05", "3.34", "0.1286", 8),
    @(2, "014668", "银华专精特新量化优选股票A", "0.26", "94.19", "1.33", "0.0035", 7),
    @(3, "003308", "中信建投睿利灵活配置混合A", "0.07", "71.26", "3.13", "0.0022", 7),
    @(4, "014669", "银华专精特新量化优选股票C", "0.15", "94.19", "1.33", "0.0020", 7),
    @(5, "004635", "中信建投睿利灵活配置混合C", "0.04", "71.26", "3.13", "0.0013", 7)
)

$row = 2
foreach ($r in $q4Data) {
    $q4.Cells.Item($row, 1).Value = $r[0]
    $q4.Cells.Item($row, 2).Value = $r[1]
    $q4.Cells.Item($row, 3).Value = $r[2]
    $q4.Cells.Item($row, 4).Value = $r[3]
    $q4.Cells.Item($row, 5).Value = $r[4]
    $q4.Cells.Item($row, 6).Value = $r[5]
    $q4.Cells.Item($row, 7).Value = $r[6]
    $q4.Cells.Item($row, 8).Value = $r[7]
    $row++
}

# ---------------------------------------------------------------------
# Step 2: update "总计" - insert the 2022-Q4 total as the new row 2, and
# shift 2022-Q3 / 2022-Q2 down to rows 3 and 4.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 6
$total.Cells.Item(2, 4).Value = 0.31

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2022-Q3"
$total.Cells.Item(3, 3).Value = 9
$total.Cells.Item(3, 4).Value = 1.32

# Row 4 is brand new - copy formatting (style) from row 3's A cell first.
$total.Range("A3").Copy()
$total.Range("A4").PasteSpecial(-4122)

$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = "2022-Q2"
$total.Cells.Item(4, 3).Value = 3
$total.Cells.Item(4, 4).Value = 0.04

# Restore the original active sheet / selection state.
$total.Activate()
